$wb = $excel.ActiveWorkbook

# --- Norway: copy the Hungary sheet (same template) and place right after it ---
$hungary = $wb.Worksheets.Item("Hungary")
$hungary.Copy([System.Reflection.Missing]::Value, $hungary)
$norway = $wb.Worksheets.Item("Hungary (2)")
$norway.Name = "Norway"

# Fill in Norway-specific values (ticket id first, then market name - matches
# shared-string insertion order of the target workbook)
$norway.Range("B4").Value = "NGC-2931/T3061"
$norway.Range("B2").Value = "Norway Market"

# Rows 3-5 should use the default (non custom) row height on the new sheet
$norway.Rows.Item(3).AutoFit()
$norway.Rows.Item(4).AutoFit()
$norway.Rows.Item(5).AutoFit()

# Column D is wider on the new sheet and no longer auto "best fit"
$norway.Columns.Item(4).ColumnWidth = 17.5

# --- Poland: copy the Norway sheet and place right after it ---
$norway.Copy([System.Reflection.Missing]::Value, $norway)
$poland = $wb.Worksheets.Item("Norway (2)")
$poland.Name = "Poland"

$poland.Range("B4").Value = "NGC-2920/T3100"
$poland.Range("B2").Value = "Poland Market"

# Make the whole-sheet selection match the freshly created-sheet look
$poland.Activate()
$poland.Cells.Select()
$norway.Activate()
$norway.Cells.Select()

# Norway ends up as the active tab, as in the target workbook
$norway.Activate()
